$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("current punchlist")

# Remove the branch value from C7 (row for "master" branch note no longer set here)
$ws.Range("C7").Clear()

# Add new release-note row for the "multipath" work
$ws.Range("A19").Value2 = "WORKING"
$ws.Range("B19").Value2 = "Feature improvement"
$ws.Range("F19").Value2 = "Add capability for multiple pathogens."
$ws.Range("C19").Value2 = "multipath"
$ws.Rows.Item(19).RowHeight = 22

# Move the active selection to C20, matching where the user would continue entry
$ws.Range("C20").Select()
